$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 276, shifting existing rows 276:289 down to 277:290
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new weekly record
$ws.Range("A276").Value = 11
$ws.Range("B276").Value = "Vega Monumental Concepción"
$ws.Range("C276").Value = "Bíobío"
$ws.Range("D276").Value = 44714
$ws.Range("E276").Value = 8
$ws.Range("F276").Value = 100112023
$ws.Range("G276").Value = "Brócoli"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 2700
$ws.Range("K276").Value = 700
$ws.Range("L276").Value = 750
$ws.Range("M276").Value = 728
$ws.Range("N276").Value = "$/unidad"
$ws.Range("O276").Value = "Región Metropolitana"
$ws.Range("P276").Value = 728
$ws.Range("Q276").Value = 1
$ws.Range("R276").Value = "Hortaliza"
